$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The paragraph holding only the "_GoBack" bookmark is the very last
# paragraph in the document body (right before the final section
# break). The paragraph right before it is an existing empty one.
# Work off paragraph indices (Range.Paragraphs on a collapsed range is
# unreliable in this host), which is robust and avoids ambiguity.
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$bmPara = $d.Paragraphs.Item($lastIndex)
$prevPara = $d.Paragraphs.Item($lastIndex - 1)

# ------------------------------------------------------------------
# New block of paragraphs to insert between the existing empty
# paragraph and the paragraph holding the bookmark: two more blank
# paragraphs, then the "Ex. 0.17" exercise write-up.
# ------------------------------------------------------------------
$newBlockXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pBdr><w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/></w:pBdr></w:pPr><w:r><w:t>Ex. 0.17</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>N=20</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Searched value=5</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Iteration 1:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>start</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>=0 end=19</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>mid</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>=9</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Iteration 2</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>start</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>=0 end=8</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>mid</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>=4</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Iteration 3</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>start</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>=5 end=8</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>mid</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>=6</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Iteration 4</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>start</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>=5 end=5</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="gramStart"/><w:r><w:t>mid</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>=5</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

# Insert a fresh paragraph right after the pre-existing empty one, then
# push the whole new block's XML into it - InsertXML happily expands a
# single target paragraph into as many paragraphs as the XML contains.
$prevPara.Range.InsertParagraphAfter()
$newIndex = $lastIndex  # the freshly inserted blank paragraph's index
$d.Paragraphs.Item($newIndex).Range.InsertXML($newBlockXml)

# ------------------------------------------------------------------
# Finally, rewrite the paragraph that used to hold only the bookmark so
# that it carries the closing commentary sentence, keeping the
# "_GoBack" bookmark anchored exactly where it was (between "before
# the " and "last").
# ------------------------------------------------------------------
$finalParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>On every iteration</w:t></w:r><w:r><w:t xml:space="preserve"> after the first</w:t></w:r><w:r><w:t xml:space="preserve"> and before the </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>last</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>the searched range is halved by readjusting the values of start or end, depending of whether the searched value may appear before or after mid</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>.</w:t></w:r></w:p>
'@

$bm = $d.Bookmarks("_GoBack")
$bmParaIndex = $d.Paragraphs.Count
$d.Paragraphs.Item($bmParaIndex).Range.InsertXML($finalParaXml)
